$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 16 (template row) down to new rows 17-19
$ws.Range("A16:M16").Copy()
$ws.Range("A17:M17").PasteSpecial(-4122)
$ws.Range("A16:M16").Copy()
$ws.Range("A18:M18").PasteSpecial(-4122)
$ws.Range("A16:M16").Copy()
$ws.Range("A19:M19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Now set final values for rows 10-19 (A,B,C..M)

$ws.Range("A10").Value2 = 8
$ws.Range("B10").Value2 = "Gaussian-Quadrature"
$ws.Range("C10").Value2 = 1.781905590457717
$ws.Range("D10").Value2 = 0.1727345983995174
$ws.Range("E10").Value2 = 1.083837254209483
$ws.Range("F10").Value2 = 1.781905590457717
$ws.Range("G10").Value2 = 0.5635907596766025
$ws.Range("H10").Value2 = 1.306867444027801
$ws.Range("I10").Value2 = 1.165253631243138
$ws.Range("J10").Value2 = 0.1727345983995174
$ws.Range("K10").Value2 = 0.6282859263045003
$ws.Range("L10").Value2 = 1.205095758381109
$ws.Range("M10").Value2 = 1.012364879669043

$ws.Range("A11").Value2 = 9
$ws.Range("B11").Value2 = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value2 = 1.059044241197765
$ws.Range("D11").Value2 = 0.9342228828008236
$ws.Range("E11").Value2 = 1.089393757715615
$ws.Range("F11").Value2 = 1.059044241197765
$ws.Range("G11").Value2 = 0.7499412156325075
$ws.Range("H11").Value2 = 1.560991105673924
$ws.Range("I11").Value2 = 0.9785707088574476
$ws.Range("J11").Value2 = 0.9342228828008236
$ws.Range("K11").Value2 = 1.01180832025822
$ws.Range("L11").Value2 = 1.035426280727992
$ws.Range("M11").Value2 = 1.062027318646347

$ws.Range("A12").Value2 = 10
$ws.Range("B12").Value2 = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value2 = 1.052516074712303
$ws.Range("D12").Value2 = 0.9358853976740306
$ws.Range("E12").Value2 = 1.089179073775055
$ws.Range("F12").Value2 = 1.052516074712303
$ws.Range("G12").Value2 = 0.7512985668229757
$ws.Range("H12").Value2 = 1.556530411707626
$ws.Range("I12").Value2 = 0.9788019802851171
$ws.Range("J12").Value2 = 0.9358853976740306
$ws.Range("K12").Value2 = 1.012532235724543
$ws.Range("L12").Value2 = 1.032524155218423
$ws.Range("M12").Value2 = 1.060701917496184

$ws.Range("A13").Value2 = 11
$ws.Range("B13").Value2 = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value2 = 1.058512636069459
$ws.Range("D13").Value2 = 0.9337738471325808
$ws.Range("E13").Value2 = 1.090817955237511
$ws.Range("F13").Value2 = 1.058512636069459
$ws.Range("G13").Value2 = 0.7497796156059828
$ws.Range("H13").Value2 = 1.552777392780684
$ws.Range("I13").Value2 = 0.9793726038546173
$ws.Range("J13").Value2 = 0.9337738471325808
$ws.Range("K13").Value2 = 1.012295901185046
$ws.Range("L13").Value2 = 1.035404268627253
$ws.Range("M13").Value2 = 1.060839008446806

$ws.Range("A14").Value2 = 12
$ws.Range("B14").Value2 = "NoRotation-tilt60deg"
$ws.Range("C14").Value2 = 0.4253800000000011
$ws.Range("D14").Value2 = 0.4870880000000012
$ws.Range("E14").Value2 = 0.9422999999999986
$ws.Range("F14").Value2 = 0.4253800000000011
$ws.Range("G14").Value2 = 0.3716520000000004
$ws.Range("H14").Value2 = 6.584716
$ws.Range("I14").Value2 = 0.7553880000000004
$ws.Range("J14").Value2 = 0.4870880000000012
$ws.Range("K14").Value2 = 0.7146939999999999
$ws.Range("L14").Value2 = 0.5700370000000005
$ws.Range("M14").Value2 = 1.594420666666667

$ws.Range("A15").Value2 = 13
$ws.Range("B15").Value2 = "Rotation-NoTilt"
$ws.Range("C15").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0.8402750000000009
$ws.Range("F15").Value2 = 0
$ws.Range("G15").Value2 = 0
$ws.Range("H15").Value2 = 11.07289999999999
$ws.Range("I15").Value2 = 0.55
$ws.Range("J15").Value2 = 0
$ws.Range("K15").Value2 = 0.4201375000000004
$ws.Range("L15").Value2 = 0.2100687500000002
$ws.Range("M15").Value2 = 2.077195833333331

$ws.Range("A16").Value2 = 14
$ws.Range("B16").Value2 = "Rotation-60detTilt"
$ws.Range("C16").Value2 = 0.4176518461440009
$ws.Range("D16").Value2 = 0.4363601430528017
$ws.Range("E16").Value2 = 0.9537747055616003
$ws.Range("F16").Value2 = 0.4176518461440009
$ws.Range("G16").Value2 = 0.3921760903167999
$ws.Range("H16").Value2 = 6.578790188543997
$ws.Range("I16").Value2 = 0.7680424218623999
$ws.Range("J16").Value2 = 0.4363601430528017
$ws.Range("K16").Value2 = 0.695067424307201
$ws.Range("L16").Value2 = 0.5563596352256009
$ws.Range("M16").Value2 = 1.5911325659136

$ws.Range("A17").Value2 = 15
$ws.Range("B17").Value2 = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value2 = 0.979829505812832
$ws.Range("D17").Value2 = 0.9904905745000722
$ws.Range("E17").Value2 = 0.9922433789158267
$ws.Range("F17").Value2 = 0.979829505812832
$ws.Range("G17").Value2 = 0.9883304894907732
$ws.Range("H17").Value2 = 0.9992185099432691
$ws.Range("I17").Value2 = 0.9942745510350426
$ws.Range("J17").Value2 = 0.9904905745000722
$ws.Range("K17").Value2 = 0.9913669767079494
$ws.Range("L17").Value2 = 0.9855982412603908
$ws.Range("M17").Value2 = 0.9907311682829693

$ws.Range("A18").Value2 = 16
$ws.Range("B18").Value2 = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value2 = 1.048538032049184
$ws.Range("D18").Value2 = 1.03761618970499
$ws.Range("E18").Value2 = 1.093227248382229
$ws.Range("F18").Value2 = 1.048538032049184
$ws.Range("G18").Value2 = 1.057872244402448
$ws.Range("H18").Value2 = 1.107667779386402
$ws.Range("I18").Value2 = 0.9238249844964177
$ws.Range("J18").Value2 = 1.03761618970499
$ws.Range("K18").Value2 = 1.06542171904361
$ws.Range("L18").Value2 = 1.056979875546397
$ws.Range("M18").Value2 = 1.044791079736945

$ws.Range("A19").Value2 = 17
$ws.Range("B19").Value2 = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value2 = 0.9927611213540172
$ws.Range("D19").Value2 = 1.400361187171606
$ws.Range("E19").Value2 = 0.8993945915700277
$ws.Range("F19").Value2 = 0.9927611213540172
$ws.Range("G19").Value2 = 1.11265427153207
$ws.Range("H19").Value2 = 0.6797203513583618
$ws.Range("I19").Value2 = 0.9250598974596325
$ws.Range("J19").Value2 = 1.400361187171606
$ws.Range("K19").Value2 = 1.149877889370817
$ws.Range("L19").Value2 = 1.071319505362417
$ws.Range("M19").Value2 = 1.001658570074286
